$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I4").Value = 2020
